$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("SV-RJ47679", "00088387-2", "1416-271278-101-0", "JUAN DAVID, REYES SALAZAR", "MASCULINO", "7071-0453", "INACTIVO", "292.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-NM9278", "00090088-6", "0614-271260-014-0", "MANUEL VICENTE, NUÑEZ GRANADOS", "MASCULINO", "6131-7133", "INACTIVO", "35.20", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-RM20831", "00196268-7", "0315-180165-001-5", "MAURICIO ERNESTO, RAMIREZ MARTINEZ", "MASCULINO", "6261-0010", "ACTIVO", "149.20", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-MM47693", "00310459-4", "0702-240253-001-0", "MANUEL ANTONIO, MEJIA", "MASCULINO", "6235-7596", "INACTIVO", "189.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-OM47682", "00356539-2", "1217-241278-104-0", "MARGARITA DE LA PAZ, OLIVARES RUIZ", "FEMENINO", "6660-6154", "ACTIVO", "560.00", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-CJ47665", "00371900-6", "0515-080248-001-9", "JUAN JUVENCIO, CASTILLO MEZQUITA", "MASCULINO", "7851-4479", "ACTIVO", "171.60", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-PR20278", "00374800-4", "0821-220252-001-2", "RAUL ALBERTO, PINEDA DIAZ", "MASCULINO", "6278-7729", "ACTIVO", "302.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-SR47630", "00424078-3", "1008-050256-101-6", "RAUL, SANCHEZ MARTINEZ", "MASCULINO", "7439-8467", "ACTIVO", "185.20", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-PA13188", "00466902-6", "1011-280273-101-0", "ANGEL VICTOR, PINO MERINO", "MASCULINO", "7923-2040", "ACTIVO", "256.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-MA8266", "00559869-2", "0614-190367-101-1", "ALICIA ORBELINA, MENDEZ BENITEZ", "MASCULINO", "7070-5856", "ACTIVO", "297.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-QA47695", "00619610-0", "0614-290358-010-0", "ANA DORA ALICIA, QUINTANILLA DE ZULETA", "FEMENINO", "6271-6439", "ACTIVO", "234.20", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-MJ30105", "00630581-4", "0614-220281-103-6", "JAIME RAUL, MARROQUIN MARROQUIN", "MASCULINO", "7976-7136", "ACTIVO", "260.20", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-UD47691", "00809544-9", "1403-071155-101-1", "DOLORES SUYAPA, UMANZOR DE RIVERA", "FEMENINO", "6225-7198", "ACTIVO", "368.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-GS47519", "00829876-2", "1010-030362-001-1", "SARA LORENA, GUZMAN SAMAYOA", "FEMENINO", "6284-6340", "ACTIVO", "384.20", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-RJ47676", "00881659-4", "0512-030470-101-7", "JACOBO, RAMOS PEREZ", "MASCULINO", "7165-1054", "INACTIVO", "184.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-PG47677", "00927645-4", "0522-200472-101-3", "GUILLERMO, POCASANGRE HUEZO", "MASCULINO", "7789-4991", "ACTIVO", "272.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-CB47683", "01062229-4", "1109-250583-102-4", "BLANCA MARIBEL, CRUZ AMAYA", "FEMENINO", "7886-6026", "ACTIVO", "335.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-FR47663", "01075769-1", "0614-240560-013-6", "ROBERTO ENRIQUE, FONG HERNANDEZ", "MASCULINO", "6282-3380", "INACTIVO", "361.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-BJ15868", "01167192-9", "1217-210658-004-3", "JOSE ATILIO, BENITEZ PARADA", "MASCULINO", "6207-4477", "ACTIVO", "308.60", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-LC6418", "01421695-6", "0614-031062-003-1", "CESAR MAURICIO, LOPEZ GUZMAN", "MASCULINO", "6305-5568", "ACTIVO", "86.60", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-MN1421", "01448290-5", "0210-010765-101-4", "NESTOR MARIO, MARROQUIN AREVALO", "MASCULINO", "7740-7501", "INACTIVO", "212.00", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-GC42266", "01564617-0", "0511-121181-102-6", "CARLOS RICARDO, GHIRINGHELLO ROSALES", "FEMENINO", "7247-1706", "INACTIVO", "214.20", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-CH47710", "01695478-8", "1324-220868-101-2", "HECTOR MANUEL, CHACON ARGUETA", "MASCULINO", "7778-1096", "ACTIVO", "82.60", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-HE47712", "01771937-7", "0614-020980-118-1", "ERICK STEFFAN, HERNANDEZ DURAN", "MASCULINO", "7695-5228", "ACTIVO", "364.60", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-MG2164", "01804572-1", "0210-250264-002-1", "GABRIEL ANTONIO, MAYORGA RAMIREZ", "MASCULINO", "7070-5723", "ACTIVO", "281.20", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-GM47595", "01890496-7", "0502-250251-001-5", "MARIO, GUARDADO", "MASCULINO", "7928-9424", "ACTIVO", "318.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-LR47673", "02054051-7", "0614-190981-137-0", "RICARDO ERNESTO, LARA CACERES", "MASCULINO", "7985-7656", "INACTIVO", "223.40", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-MH23670", "02313750-9", "1114-260170-101-5", "HECTOR ANTONIO, MARTINEZ ORTIZ", "MASCULINO", "6631-1175", "ACTIVO", "141.80", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-DF22100", "02377745-6", "0614-180648-005-2", "FRANCISCO ENRIQUE, DIAZ RODRIGUEZ", "MASCULINO", "6523-6600", "ACTIVO", "210.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-SM1054", "02408291-9", "0203-051068-101-4", "MARITZA CAROLINA, SALAZAR DE SERMEÑO", "FEMENINO", "6447-8435", "ACTIVO", "733.60", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-MJ25233", "02435337-8", "0614-161264-013-0", "JORGE ALBERTO, MIRANDA QUINTANILLA", "MASCULINO", "6272-9067", "INACTIVO", "214.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-ME47627", "02454429-6", "0607-270265-001-7", "EVELYN DE LOS ANGELES, MARROQUIN AMAYA", "FEMENINO", "7070-4632", "ACTIVO", "395.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-GE25309", "02575654-5", "1407-201281-101-6", "EDWIN JOVANY, GONZALEZ", "MASCULINO", "7410-0034", "ACTIVO", "163.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-CR47598", "02593311-3", "0614-200854-003-2", "ROSA MARGARITA, CRUZ GUZMAN", "FEMENINO", "6243-7362", "ACTIVO", "363.40", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-PJ19861", "02726756-8", "0210-170964-007-2", "JOSE MAURICIO, PINEDA MULATO", "MASCULINO", "6232-1846", "ACTIVO", "369.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-OL19455", "02921225-0", "0103-240966-101-0", "LUIS WALTER, ORTIZ MEDINA", "MASCULINO", "6294-7980", "ACTIVO", "251.80", "CLIENTE EJECUTIVO", "sv")
    ,@("SV-RR5909", "03365212-9", "0610-150982-107-3", "ROBERTO, RAMIREZ DEODANES", "MASCULINO", "7536-9733", "ACTIVO", "143.60", "CLIENTE PREFERENCIAL", "sv")
    ,@("SV-MS47711", "04326581-9", "1123-011290-101-6", "SAMUEL GERARDO, MEDRANO RODRIGUEZ", "MASCULINO", "7869-2068", "ACTIVO", "219.40", "CLIENTE EJECUTIVO", "sv")
)

$startRow = 78
$endRow = $startRow + $data.Count - 1

# Column H holds decimal-looking text (e.g. "292.40"); format as Text first
# so Excel does not silently coerce the literal into a Number.
$hRange = $ws.Range($ws.Cells.Item($startRow, 8), $ws.Cells.Item($endRow, 8))
$hRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Restore the plain/Normal cell style now that the text values are locked in,
# so the new rows do not carry a lingering custom style index.
$hRange.Style = "Normal"

Write-Host "Done writing rows $startRow to $endRow"
